$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "log over-top10"
$ws.Range("H2").Value = 0.94299999999999995
$ws.Range("H3").Value = 0.050500000000000003
$ws.Range("H4").Value = 0.025999999999999999
$ws.Range("H5").Value = 0.871
$ws.Range("H6").Value = 0.90710000000000002

$ws.Range("H2:H6").NumberFormat = "0.00%"

$ws.Range("H7").Select()
